$wb = $excel.ActiveWorkbook

# --- Sheet "Transaksi" (first sheet) ---
$ws1 = $wb.Worksheets.Item(1)

# Row 2
$ws1.Range("A2").Value = 45443
$ws1.Range("A2").NumberFormat = "YYYY-MM-DD"
$ws1.Range("D2").Value = "148.943.002,00"
$ws1.Range("E2").Value = 45443
$ws1.Range("E2").NumberFormat = "YYYY-MM-DD"
$ws1.Range("H2").Value = "148.943.002,00"

# Row 3
$ws1.Range("A3").Value = 45447
$ws1.Range("A3").NumberFormat = "YYYY-MM-DD"
$ws1.Range("B3").Value = "114.574.866,00"
$ws1.Range("D3").Value = "263.517.868,00"
$ws1.Range("E3").Value = 45447
$ws1.Range("E3").NumberFormat = "YYYY-MM-DD"
$ws1.Range("G3").Value = "114.574.866,00"
$ws1.Range("H3").Value = "263.517.868,00"

# Row 4
$ws1.Range("A4").Value = 45473
$ws1.Range("A4").NumberFormat = "YYYY-MM-DD"
$ws1.Range("C4").Value = "36.709,00"
$ws1.Range("D4").Value = "263.481.159,00"
$ws1.Range("E4").Value = 45473
$ws1.Range("E4").NumberFormat = "YYYY-MM-DD"
$ws1.Range("F4").Value = "36.709,00"
$ws1.Range("H4").Value = "263.481.159,00"
$ws1.Range("J4").Value = "0,00"
$ws1.Range("K4").Value = "Matched"
$ws1.Range("L4").Value = "-"

# Row 5
$ws1.Range("A5").Value = 45473
$ws1.Range("A5").NumberFormat = "YYYY-MM-DD"
$ws1.Range("C5").Value = "12.000,00"
$ws1.Range("D5").Value = "263.469.159,00"
$ws1.Range("E5").Value = 45473
$ws1.Range("E5").NumberFormat = "YYYY-MM-DD"
$ws1.Range("F5").Value = "12.000,00"
$ws1.Range("H5").Value = "263.469.159,00"
$ws1.Range("J5").Value = "0,00"
$ws1.Range("K5").Value = "Matched"
$ws1.Range("L5").Value = "-"

# Row 6
$ws1.Range("A6").Value = 45473
$ws1.Range("A6").NumberFormat = "YYYY-MM-DD"
$ws1.Range("C6").Value = "10.000,00"
$ws1.Range("D6").Value = "263.459.159,00"
$ws1.Range("E6").Value = 45473
$ws1.Range("E6").NumberFormat = "YYYY-MM-DD"
$ws1.Range("F6").Value = "10.000,00"
$ws1.Range("H6").Value = "263.459.159,00"
$ws1.Range("J6").Value = "0,00"
$ws1.Range("K6").Value = "Matched"
$ws1.Range("L6").Value = "-"

# Row 7
$ws1.Range("A7").Value = 45473
$ws1.Range("A7").NumberFormat = "YYYY-MM-DD"
$ws1.Range("B7").Value = "183.543,00"
$ws1.Range("D7").Value = "263.642.702,00"
$ws1.Range("E7").Value = 45473
$ws1.Range("E7").NumberFormat = "YYYY-MM-DD"
$ws1.Range("G7").Value = "183.543,00"
$ws1.Range("H7").Value = "263.642.702,00"

# Row 8
$ws1.Range("A8").Value = 45473
$ws1.Range("A8").NumberFormat = "YYYY-MM-DD"
$ws1.Range("B8").Value = "114.758.409,00"
$ws1.Range("C8").Value = "58.709,00"
$ws1.Range("D8").Value = "263.642.702,00"
$ws1.Range("E8").Value = 45473
$ws1.Range("E8").NumberFormat = "YYYY-MM-DD"
$ws1.Range("F8").Value = "58.709,00"
$ws1.Range("G8").Value = "114.758.409,00"
$ws1.Range("H8").Value = "263.642.702,00"
$ws1.Range("K8").Value = "Closing Balance"
$ws1.Range("L8").Value = ""

# Remove the old rows 9, 10 and 11 (data now ends at row 8)
$ws1.Rows.Item(9).EntireRow.Delete()
$ws1.Rows.Item(9).EntireRow.Delete()
$ws1.Rows.Item(9).EntireRow.Delete()

# --- Sheet "Summary" (second sheet) ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B2").Value = 45443
$ws2.Range("B2").NumberFormat = "YYYY-MM-DD"
$ws2.Range("C2").Value = "148.943.002,00"
$ws2.Range("D2").Value = "148.943.002,00"

$ws2.Range("B3").Value = 45473
$ws2.Range("B3").NumberFormat = "YYYY-MM-DD"
$ws2.Range("C3").Value = "263.642.702,00"
$ws2.Range("D3").Value = "263.642.702,00"
